# NatmiData LR-pairs TPM update: Sending cluster "MuSCs" -> "Resolving-Mac",
# refreshed TPM-derived metrics, and a new self-pairing row (Resolving-Mac -> Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sending cluster" column: MuSCs -> Resolving-Mac for all existing data rows.
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("A4").Value = "Resolving-Mac"

# Row 2 (Target cluster = ECs): refreshed TPM-derived values.
$ws.Range("G2").Value = 0.02530666666666667
$ws.Range("H2").Value = 0.07592
$ws.Range("M2").Value = 4.618552666666667
$ws.Range("N2").Value = 13.855658
$ws.Range("O2").Value = 0.78434648953826
$ws.Range("P2").Value = 0.78434648953826
$ws.Range("Q2").Value = 0.1168801728177778
$ws.Range("R2").Value = 1.05192155536
$ws.Range("S2").Value = 0.78434648953826
$ws.Range("T2").Value = 0.78434648953826

# Row 3 (Target cluster = FAPs): refreshed TPM-derived values.
$ws.Range("G3").Value = 0.02530666666666667
$ws.Range("H3").Value = 0.07592
$ws.Range("O3").Value = 0.1153620112191035
$ws.Range("P3").Value = 0.1153620112191036
$ws.Range("Q3").Value = 0.01719078492444445
$ws.Range("R3").Value = 0.15471706432
$ws.Range("S3").Value = 0.1153620112191035
$ws.Range("T3").Value = 0.1153620112191036

# Row 4 (Target cluster stays self-pairing, now named "MuSCs" again since the
# sending cluster took the "Resolving-Mac" name): refreshed TPM-derived values.
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.02530666666666667
$ws.Range("H4").Value = 0.07592
$ws.Range("M4").Value = 0.478937
$ws.Range("N4").Value = 1.436811
$ws.Range("O4").Value = 0.0813355572127976
$ws.Range("P4").Value = 0.08133555721279762
$ws.Range("Q4").Value = 0.01212029901333333
$ws.Range("R4").Value = 0.10908269112
$ws.Range("S4").Value = 0.0813355572127976
$ws.Range("T4").Value = 0.08133555721279762

# New row 5: Resolving-Mac -> Resolving-Mac self-pairing interaction.
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Fgf8"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02530666666666667
$ws.Range("H5").Value = 0.07592
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1116203333333333
$ws.Range("N5").Value = 0.334861
$ws.Range("O5").Value = 0.01895594202983873
$ws.Range("P5").Value = 0.01895594202983874
$ws.Range("Q5").Value = 0.002824738568888889
$ws.Range("R5").Value = 0.02542264712
$ws.Range("S5").Value = 0.01895594202983873
$ws.Range("T5").Value = 0.01895594202983874
